$d = $word.ActiveDocument

$replacements = @(
    @{old="315×5=1575"; new="284×2=568"},
    @{old="218×7=1526"; new="577×2=1154"},
    @{old="970×2=1940"; new="179×2=358"},
    @{old="767×2=1534"; new="448×2=896"},
    @{old="876×9=7884"; new="216×3=648"},
    @{old="724×4=2896"; new="910×9=8190"},
    @{old="976×3=2928"; new="192×8=1536"},
    @{old="323×4=1292"; new="484×6=2904"},
    @{old="469×2=938"; new="443×4=1772"},
    @{old="745×7=5215"; new="424×2=848"},
    @{old="460×8=3680"; new="862×3=2586"},
    @{old="468×3=1404"; new="274×3=822"},
    @{old="721×4=2884"; new="785×4=3140"},
    @{old="317×6=1902"; new="268×5=1340"},
    @{old="631×7=4417"; new="102×3=306"},
    @{old="956×2=1912"; new="895×7=6265"},
    @{old="953×6=5718"; new="913×2=1826"},
    @{old="161×8=1288"; new="129×4=516"},
    @{old="267×7=1869"; new="941×5=4705"},
    @{old="498×8=3984"; new="705×5=3525"},
    @{old="229×4=916"; new="390×6=2340"},
    @{old="744×7=5208"; new="375×2=750"},
    @{old="437×4=1748"; new="520×2=1040"},
    @{old="435×9=3915"; new="479×6=2874"},
    @{old="837×9=7533"; new="331×3=993"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
